$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "Scaling factor" column (R) ---------------------------------------
$ws.Range("R1").Value = "Scaling factor"
$ws.Range("R1").Font.Bold = $true

# Scaling factors per building row
$ws.Range("R2").Value = 0.3819
$ws.Range("R3").Value = 0.5875
$ws.Range("R4").Value = 0.7632
$ws.Range("R2:R4").Font.Bold = $true

# --- Convert the constant columns F, G, H, I, J, M, N, O into formulas -----
# that scale the original constant by the new per-row scaling factor.
$ws.Range("F2").Formula = "=174.14*R2"
$ws.Range("G2").Formula = "=286.528*R2"
$ws.Range("H2").Formula = "=115.621*R2"
$ws.Range("I2").Formula = "=49.261*R2"
$ws.Range("J2").Formula = "=213505.516*R2"
$ws.Range("M2").Formula = "=10.037*R2"
$ws.Range("N2").Formula = "=4.533*R2"
$ws.Range("O2").Formula = "=1.619*R2"

$ws.Range("F3").Formula = "=174.14*R3"
$ws.Range("G3").Formula = "=286.528*R3"
$ws.Range("H3").Formula = "=115.621*R3"
$ws.Range("I3").Formula = "=49.261*R3"
$ws.Range("J3").Formula = "=213505.516*R3"
$ws.Range("M3").Formula = "=10.037*R3"
$ws.Range("N3").Formula = "=4.533*R3"
$ws.Range("O3").Formula = "=1.619*R3"

$ws.Range("F4").Formula = "=174.14*R4"
$ws.Range("G4").Formula = "=286.528*R4"
$ws.Range("H4").Formula = "=115.621*R4"
$ws.Range("I4").Formula = "=49.261*R4"
$ws.Range("J4").Formula = "=213505.516*R4"
$ws.Range("M4").Formula = "=10.037*R4"
$ws.Range("N4").Formula = "=4.533*R4"
$ws.Range("O4").Formula = "=1.619*R4"

# --- View / window cosmetic changes -----------------------------------------
$win = $excel.ActiveWindow
$win.Zoom = 150
$win.ScrollRow = 1
$win.ScrollColumn = 11
$ws.Range("P10").Select()

$wb.Save()
